$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-04-06"

# Update the 2022 column header label (I1)
$ws.Range("I1").Value = "2022 (through 04-06)"

# Add data for 2022-04-14: bump May (row 5) 2022 value by 1
$ws.Range("I5").Value = 17

# Update the running Total (row 14) 2022 value by 1
$ws.Range("I14").Value = 451
